$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete row 2 (the sub-header row containing "(m3/s)", "(MW)", "(MW)", "(GWh)", "(GWh)", "(GWh)")
# This shifts all data rows (3..104) up by one, so the plant data now starts at row 2.
$ws.Rows.Item(2).Delete()

# Update the selection to match the target state: active cell A2, selection A2:K2
$ws.Range("A2:K2").Select()
